# Append/update "2026-01-04 06:28 JST" scrape results into the "ランサーズ" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Column width tweaks ---------------------------------------------------
# ColumnWidth is specified in "characters" (Normal-style digit width) while
# the OOXML <col width="..."> attribute that ends up on disk is that value
# re-expressed in the MDW-7 pixel grid (stored = chars + 11/12, quantized to
# 1/256ths). Feeding in target-1+1/6 round-trips to an exact integer stored
# width (47 / 30 / 12) instead of the fractional ...833333 you'd get by
# naively assigning the integer target itself.
$ws.Columns.Item(2).ColumnWidth = 46 + 1/6   # -> stored width 47
$ws.Columns.Item(4).ColumnWidth = 29 + 1/6   # -> stored width 30
$ws.Columns.Item(8).ColumnWidth = 11 + 1/6   # -> stored width 12

# --- Drop all existing hyperlinks up front (this host clears the whole
#     sheet's collection no matter which range it's invoked on, so do it
#     once before anything else needs its own hyperlink) ------------------
$ws.Range("F2:F7").Hyperlinks.Delete()

# --- Row 2 ----------------------------------------------------------------
$ws.Range("A2").Value = "2026-01-04 06:28:45"
$ws.Range("B2").Value = "【急募】医療診断AIの深層学習モデル開発の専門家募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5464587"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Row 3 ----------------------------------------------------------------
$ws.Range("A3").Value = "2026-01-04 06:28:45"
$ws.Range("B3").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G3").Value = 243
$ws.Range("H3").Value = "🔥API ◆ツール"

# --- Row 4 ----------------------------------------------------------------
$ws.Range("A4").Value = "2026-01-04 06:28:45"
$ws.Range("B4").Value = "【急募】簡単なHP作成とAWS構築をしてくれる方募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5457524"
$ws.Range("G4").Value = 18
$ws.Range("H4").Clear()

# --- Drop the now-stale rows 5:7 -------------------------------------------
$ws.Range("A5:A7").EntireRow.Delete()

# --- Re-create hyperlinks for the surviving URL cells ----------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5464587")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5457524")

# Hyperlinks.Add stamps a brand-new (duplicate) cell style; pin these cells
# back onto the workbook's single shared "Hyperlink" cell style so we don't
# leave redundant style records referenced in the sheet.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
